$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-22 02:38:50"
$wsZhCn.Range("H4").Value = "2016-03-22 02:39:19"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-22 02:38:53"
$wsDeDe.Range("H4").Value = "2016-03-22 02:39:25"
